$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price/volume refresh). Values are written with a
# leading apostrophe so numeric-looking strings (e.g. "0.569") stay text
# instead of being parsed into doubles, then the cell Style is reset to
# "Normal" so no stray quote-prefix style index is introduced - matching
# the source cells, which carry no explicit style.
$ws.Range("D2").Value = "'51.974.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.87%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.933.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'357.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.53%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'110.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.68%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.41%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.633"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.98%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.91%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0880"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.07%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.89%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'19.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.12%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.395.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.51%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.936.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.59%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'51.999.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.85%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.33%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.61%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0984"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'270.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.33%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.76%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.189"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +14.66%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'27.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.62%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +16.74%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +13.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'10.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'38.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.39%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.11%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'52.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.46%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.85%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E38").Value = "'  -2.24%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'18.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.85%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +1.12%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +2.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'23.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'119.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.87%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.28%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.84%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.05%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.138.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.64%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -7.81%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'9.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.43%  "
$ws.Range("E51").Style = "Normal"
